$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.897.24"
$ws.Range('E2').Value = '  +5.08%  '
$ws.Range('D3').Value = "'3.846.71"
$ws.Range('E3').Value = '  +6.80%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').Value = "'423.55"
$ws.Range('E5').Value = '  +5.12%  '
$ws.Range('D6').Value = "'129.80"
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('D7').Value = "'3.836.77"
$ws.Range('E7').Value = '  +6.67%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').Value = "'0.721"
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = "'0.0000336"
$ws.Range('E12').Value = '  +10.22%  '
$ws.Range('D13').Value = "'40.68"
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').Value = "'10.26"
$ws.Range('E14').Value = '  +4.91%  '
$ws.Range('D15').Value = "'4.456.29"
$ws.Range('E15').Value = '  +6.02%  '
$ws.Range('D16').Value = "'15.78"
$ws.Range('E16').Value = '  +19.27%  '
$ws.Range('D17').Value = "'3.850.72"
$ws.Range('E17').Value = '  +7.39%  '
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = "'19.82"
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = "'67.202.11"
$ws.Range('E20').Value = '  +5.14%  '
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').Value = "'415.13"
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = "'14.92"
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').Value = "'84.11"
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('D26').Value = "'37.46"
$ws.Range('E26').Value = '  +6.55%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = "'3.23"
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').Value = "'9.85"
$ws.Range('E28').Value = '  +5.63%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = "'5.32"
$ws.Range('E29').Value = '  +3.64%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = "'9.24"
$ws.Range('E30').Value = '  +34.82%  '
$ws.Range('D31').Value = "'733.32"
$ws.Range('E31').Value = '  +9.06%  '
$ws.Range('D32').Value = "'13.08"
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('E33').Value = '  +6.38%  '
$ws.Range('E34').Value = '  +2.90%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('E36').Value = '  -4.25%  '
$ws.Range('E37').Value = '  -4.74%  '
$ws.Range('E38').Value = '  +26.06%  '
$ws.Range('D39').Value = "'55.44"
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0460"
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = "'0.0₃0728"
$ws.Range('E41').Value = '  +16.02%  '
$ws.Range('D42').Value = "'2.87"
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('E45').Value = '  -3.11%  '
$ws.Range('E46').Value = '  +9.98%  '
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = "'140.74"
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('E51').Value = '  +0.90%  '
